$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row appended to the ledger (row 6)
$ws.Range("A6").Value = 42649.64472222222
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 10000.67
$ws.Range("D6").Value = 10014.69
$ws.Range("E6").Value = 77.349997999999999
$ws.Range("F6").Value = 77.569999999999993
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.28000000000000003
$ws.Range("I6").Value = $false

# Match the date/time display used by the other rows in column A (style s="1")
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"
